# Update the date heading and the table of division problems/answers.
$d = $word.ActiveDocument

$replacements = @(
    @("2024-10-28 Monday", "2024-10-29 Tuesday"),
    @("94÷3=31, 1", "93÷5=18, 3"),
    @("47÷8=5, 7", "85÷8=10, 5"),
    @("21÷7=3, 0", "31÷4=7, 3"),
    @("34÷8=4, 2", "52÷5=10, 2"),
    @("32÷4=8, 0", "57÷8=7, 1"),
    @("37÷8=4, 5", "31÷7=4, 3"),
    @("80÷9=8, 8", "25÷6=4, 1"),
    @("23÷4=5, 3", "36÷4=9, 0"),
    @("44÷9=4, 8", "87÷6=14, 3"),
    @("85÷5=17, 0", "28÷2=14, 0"),
    @("35÷6=5, 5", "37÷7=5, 2"),
    @("80÷6=13, 2", "70÷7=10, 0"),
    @("97÷3=32, 1", "85÷3=28, 1"),
    @("30÷3=10, 0", "60÷5=12, 0"),
    @("28÷2=14, 0", "65÷7=9, 2"),
    @("98÷9=10, 8", "86÷2=43, 0"),
    @("51÷6=8, 3", "79÷6=13, 1"),
    @("48÷3=16, 0", "12÷3=4, 0"),
    @("76÷8=9, 4", "58÷8=7, 2"),
    @("53÷5=10, 3", "37÷6=6, 1"),
    @("84÷8=10, 4", "81÷2=40, 1"),
    @("50÷4=12, 2", "34÷8=4, 2"),
    @("24÷3=8, 0", "61÷3=20, 1"),
    @("40÷4=10, 0", "37÷9=4, 1"),
    @("53÷7=7, 4", "77÷7=11, 0")
)

# Two-phase replace: first swap each "old" value for a unique placeholder
# token, then swap every placeholder for its final "new" value. This avoids
# accidental re-matching when one entry's new text happens to equal another
# entry's old text (a chained collision), since Find re-scans the whole
# document on every call.
$i = 0
foreach ($pair in $replacements) {
    $old = $pair[0]
    $placeholder = "@@PLACEHOLDER_$i@@"
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $placeholder, 2)
    $i = $i + 1
}

$i = 0
foreach ($pair in $replacements) {
    $new = $pair[1]
    $placeholder = "@@PLACEHOLDER_$i@@"
    $range = $d.Content
    $range.Find.Execute($placeholder, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    $i = $i + 1
}

$d.Save()
